# "update tag add interface"
# Adds a `createUser` (User) field to GetActivity's participant info and to
# GetPost's Answer structure (which itself is simplified to reference the
# shared User type instead of repeating userId/userName/userPhoto inline).
# Also marks several interfaces as "done" on the Index sheet and refreshes
# various sheet selections / the active tab, matching the authored edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Index sheet: mark CreateActivity / GetActivity / CreatePost as "done"
# ---------------------------------------------------------------------
$wsIndex = $wb.Sheets.Item("Index")
$wsIndex.Range("C6").Value = "done"
$wsIndex.Range("C8").Value = "done"
$wsIndex.Range("C9").Value = "done"
$wsIndex.Range("C9").Select()

# ---------------------------------------------------------------------
# CreateActivity: just a cursor move, no longer the active tab
# ---------------------------------------------------------------------
$wsCreateActivity = $wb.Sheets.Item("CreateActivity")
$wsCreateActivity.Range("A13").Select()

# ---------------------------------------------------------------------
# CreatePost: cursor move only
# ---------------------------------------------------------------------
$wsCreatePost = $wb.Sheets.Item("CreatePost")
$wsCreatePost.Range("A8").Select()

# ---------------------------------------------------------------------
# GetActivity: add "picture"/String/F and "createUser"/User/T fields
# right before the nested User object, pushing it (and its members) down
# two rows.
# ---------------------------------------------------------------------
$wsGetActivity = $wb.Sheets.Item("GetActivity")

# Snapshot the nested User block before it gets overwritten.
$gaUserHeader = $wsGetActivity.Range("A20").Value

$wsGetActivity.Range("A25").Value = $wsGetActivity.Range("A23").Value
$wsGetActivity.Range("C25").Value = $wsGetActivity.Range("C23").Value
$wsGetActivity.Range("E25").Value = $wsGetActivity.Range("E23").Value

$wsGetActivity.Range("A24").Value = $wsGetActivity.Range("A22").Value
$wsGetActivity.Range("C24").Value = $wsGetActivity.Range("C22").Value
$wsGetActivity.Range("E24").Value = $wsGetActivity.Range("E22").Value

$wsGetActivity.Range("A23").Value = $wsGetActivity.Range("A21").Value
$wsGetActivity.Range("C23").Value = $wsGetActivity.Range("C21").Value
$wsGetActivity.Range("E23").Value = $wsGetActivity.Range("E21").Value

$wsGetActivity.Range("A21:F21").ClearContents()

$wsGetActivity.Range("A22").Value = $gaUserHeader

$wsGetActivity.Range("A19").Value = "picture"
$wsGetActivity.Range("C19").Value = "String"
$wsGetActivity.Range("E19").Value = "F"

$wsGetActivity.Range("A20").Value = "createUser"
$wsGetActivity.Range("C20").Value = "User"
$wsGetActivity.Range("E20").Value = "T"

$wsGetActivity.Rows(20).Select()

# ---------------------------------------------------------------------
# GetPost: add "createUser"/User/T before the Answer object, and collapse
# Answer's userId/userName/userPhoto trio into a single "user"/User/T
# reference field.
# ---------------------------------------------------------------------
$wsGetPost = $wb.Sheets.Item("GetPost")

$wsGetPost.Range("A18:F18").ClearContents()
$wsGetPost.Range("A19").Value = "Answer"

$wsGetPost.Range("A20").Value = "user"
$wsGetPost.Range("C20").Value = "User"
$wsGetPost.Range("E20").Value = "T"

$wsGetPost.Range("A21").Value = "photoList"
$wsGetPost.Range("C21").Value = "List<String(Base64)>"
$wsGetPost.Range("E21").Value = "F"

$wsGetPost.Range("A22").Value = "answer"
$wsGetPost.Range("C22").Value = "String"
$wsGetPost.Range("E22").Value = "T"

$wsGetPost.Range("A23").Value = "answerTime"
$wsGetPost.Range("C23").Value = "DateTime"
$wsGetPost.Range("E23").Value = "T"

$wsGetPost.Range("A24:F24").ClearContents()

$wsGetPost.Range("A17").Value = "createUser"
$wsGetPost.Range("C17").Value = "User"
$wsGetPost.Range("E17").Value = "T"

$wsGetPost.Activate()
$wsGetPost.Range("C30").Select()
